$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Neutro / Hc / C5ar2 -> ECs) ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1527673333333333
$ws.Range("H2").Value = 0.458302
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.002486
$ws.Range("N2").Value = 0.007458
$ws.Range("O2").Value = 0.0000300751452789208
$ws.Range("P2").Value = 0.0000300751452789208
$ws.Range("Q2").Value = 0.0003797795906666667
$ws.Range("R2").Value = 0.003418016316
$ws.Range("S2").Value = 0.0000300751452789208
$ws.Range("T2").Value = 0.0000300751452789208

# --- Row 3 (Neutro / Hc / C5ar2 -> FAPs) ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1527673333333333
$ws.Range("H3").Value = 0.458302
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01171
$ws.Range("N3").Value = 0.03513
$ws.Range("O3").Value = 0.0001416653062011917
$ws.Range("P3").Value = 0.0001416653062011917
$ws.Range("Q3").Value = 0.001788905473333333
$ws.Range("R3").Value = 0.01610014926
$ws.Range("S3").Value = 0.0001416653062011917
$ws.Range("T3").Value = 0.0001416653062011917

# --- Row 4 (Neutro / Hc / C5ar2 -> M1), replacing old row 4 (-> Neutro) ---
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1527673333333333
$ws.Range("H4").Value = 0.458302
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.107756666666667
$ws.Range("N4").Value = 27.32327
$ws.Range("O4").Value = 0.1101838716472484
$ws.Range("P4").Value = 0.1101838716472484
$ws.Range("Q4").Value = 1.391367698615556
$ws.Range("R4").Value = 12.52230928754
$ws.Range("S4").Value = 0.1101838716472484
$ws.Range("T4").Value = 0.1101838716472484

# --- Row 5 (new): Neutro / Hc / C5ar2 -> M2 ---
$ws.Range("A5").Value = "Neutro"
$ws.Range("B5").Value = "Hc"
$ws.Range("C5").Value = "C5ar2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1527673333333333
$ws.Range("H5").Value = 0.458302
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.275423333333334
$ws.Range("N5").Value = 24.82627
$ws.Range("O5").Value = 0.1001144645995861
$ws.Range("P5").Value = 0.1001144645995861
$ws.Range("Q5").Value = 1.264214354837778
$ws.Range("R5").Value = 11.37792919354
$ws.Range("S5").Value = 0.1001144645995861
$ws.Range("T5").Value = 0.1001144645995861

# --- Row 6 (new): Neutro / Hc / C5ar2 -> Neutro ---
$ws.Range("A6").Value = "Neutro"
$ws.Range("B6").Value = "Hc"
$ws.Range("C6").Value = "C5ar2"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1527673333333333
$ws.Range("H6").Value = 0.458302
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 65.26224133333334
$ws.Range("N6").Value = 195.786724
$ws.Range("O6").Value = 0.7895299233016855
$ws.Range("P6").Value = 0.7895299233016855
$ws.Range("Q6").Value = 9.969938575849778
$ws.Range("R6").Value = 89.729447182648
$ws.Range("S6").Value = 0.7895299233016855
$ws.Range("T6").Value = 0.7895299233016855
